# Apply "Add files via upload" edit: fix typo in A2, add new KB rows 5-10,
# apply thin-border / wrap-text styling, custom row heights, and update selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    'A1' = 'Issue Title'
    'B1' = 'Description'
    'C1' = 'Resolution'
    'D1' = 'Tags'
    'A2' = 'Escalation Access'
    'B2' = 'User need Escalation approval where Approve/Reject button is visible in case approval history'
    'C2' = 'We need to add the user in the Case Escalation Request Queue'
    'D2' = 'Request access'
    'A3' = 'Credit Hold'
    'B3' = 'Customer not able to create cases from SC'
    'C3' = 'Send and email to DLcollections@csod.com.'
    'D3' = 'Question'
    'A4' = 'Not able to add DSC'
    'B4' = 'Getting this error "This contact has the same email address as an active Support Central user. Duplicate Support Central users cannot be created, so this contact cannot be added as a Designated Support Contact. Please contact IT Support for assistance."'
    'C4' = 'Search the contact from User in salesforce and add .Invalid/Test at the end of the email and save'
    'D4' = 'Issue'
    'A5' = 'Portal deactivation request'
    'B5' = 'Caller wants to deactivate the portal'
    'C5' = 'This it not an official part of our Support Central management process. Normally we wait until credit hold kick in (which comes from Finance team) and then we block Support Central as part of this process.  INTERNAL:Suspending Support for Accounts on Credit Hold
If there is a reason we should step outside the normal process please escalate to Lucy Bonilla.
You could reach out to DLcollections@csod.com if you want to put them on credit hold right away since they are no longer paying their bills, that would kick off this process that blocks their access to support central with the warning shown in this KA  INTERNAL:Suspending Support for Accounts on Credit Hold 
'
    'D5' = 'Question'
    'A6' = 'JIT Authorization error'
    'C6' = 'Refere to this KA: CSOD SAGA SSO Tool throwing JIT Authorization Error'
    'D6' = 'Issue'
    'A7' = 'Not receiving correct amount of Cornerstone to Salesforce 2.0 licenses'
    'B7' = 'Cornerstone to Salesforce 2.0 licenses'
    'C7' = 'We won''t have a better process than "just ping Matt Idell". Since this is managed in a legacy Salesforce org, we can''t build out a process for this. Please refer INC0157545.'
    'D7' = 'Question'
    'A8' = 'Product Specific issues'
    'B8' = 'Product Specific issues'
    'C8' = 'This looks to be a product issue, from my time in supporting the LMS product I know that certifications is done in the product, not in Salesforce. In Support ops we can help with Salesforce and salesforce-connected apps, but we do not have access to or knowledge of the various product platforms.
If you or your customer needs product support, you can open a Support case from within Salesforce with these instructions:
1) Open case with yourself as the Contact to receive all notifications and optionally choose to hide this case from customer by unchecking ''visible in Support central'' checkbox as described in this article:https://cornerstoneondemand.lightning.force.com/lightning/articles/Knowledge/INTERNAL-How-to-Open-a-Customer-Case-With-CSOD-Employee-as-Contact
2) Open case with the customer as the contact so they receive all the notifications:https://cornerstoneondemand.lightning.force.com/lightning/articles/Knowledge/How-to-Create-a-Support-Case-for-a-Customer
I''ve set this to Resolved but please reply back if you need anything else.'
    'D8' = 'Question'
    'A9' = 'When system is comletely down'
    'C9' = 'emailing servicedeskescalations@csod.com will create a ServiceNow ticket when Okta is not available  also calling 1-800-516-8127 and leaving a voicemail will create a SNOW ticket for IT'
    'D9' = 'Issue'
    'A10' = 'Pen Test'
    'B10' = 'Pen Test request'
    'C10' = 'Refer to this page: https://cornerstoneondemandinc.service-now.com/csc?id=kb_article&table=kb_knowledge&sys_kb_id=afc7301a93c592944e0bfe9a7bba102b&spa=1'
    'D10' = 'Question'
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

# Thin box border around cells that mark a new KB-article block
$borderRanges = 'A5','A8','B8','A9','C9','A10'
foreach ($addr in $borderRanges) {
    $rng = $ws.Range($addr)
    $rng.Borders.LineStyle = 1
    $rng.Borders.Weight = 2
}

# Wrap text for the long description/resolution cells
$wrapRanges = 'C5','C8'
foreach ($addr in $wrapRanges) {
    $ws.Range($addr).WrapText = $true
}

# Row heights sized for the wrapped content
$ws.Rows.Item(5).RowHeight = 41.5
$ws.Rows.Item(8).RowHeight = 16

# Final selection left where the next new row would start
$ws.Range('D11').Select()
